$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 10.72688076740771
$ws.Range("C2").Value = 3.963739053596521
$ws.Range("D2").Value = 8.828761819488472
$ws.Range("E2").Value = 13.59623682752056
$ws.Range("F2").Value = 34.83378678937876
$ws.Range("H2").Value = 7.344005520526261
$ws.Range("J2").Value = 10.19193600088257
$ws.Range("K2").Value = 10.14184448770893
$ws.Range("M2").Value = 15.39950898954437
$ws.Range("N2").Value = 20.24264389651784
$ws.Range("O2").Value = 26.46393044682836

# Row 3
$ws.Range("B3").Value = 10.50258162370887
$ws.Range("C3").Value = 3.790300172835785
$ws.Range("D3").Value = 8.79936673783137
$ws.Range("E3").Value = 13.59191496506115
$ws.Range("F3").Value = 34.88930835545298
$ws.Range("H3").Value = 7.344005520526261
$ws.Range("J3").Value = 10.21128718076302
$ws.Range("K3").Value = 9.996660893980666
$ws.Range("M3").Value = 15.34672583660969
$ws.Range("N3").Value = 20.30316119320696
$ws.Range("O3").Value = 26.53363303243371

# Row 4
$ws.Range("B4").Value = 10.3645425397547
$ws.Range("C4").Value = 3.678900705130867
$ws.Range("D4").Value = 8.782682844002625
$ws.Range("E4").Value = 13.59151615728081
$ws.Range("F4").Value = 34.93059362212973
$ws.Range("H4").Value = 7.344005520526261
$ws.Range("J4").Value = 10.22429469031929
$ws.Range("K4").Value = 9.908205394185952
$ws.Range("M4").Value = 15.31668049222618
$ws.Range("N4").Value = 20.34205809608645
$ws.Range("O4").Value = 26.58128892995666

# Row 5
$ws.Range("B5").Value = 10.30829268747912
$ws.Range("C5").Value = 3.632307996113129
$ws.Range("D5").Value = 8.776232209314252
$ws.Range("E5").Value = 13.5919223118836
$ws.Range("F5").Value = 34.94922501810515
$ws.Range("H5").Value = 7.344005520526261
$ws.Range("J5").Value = 10.22987876008328
$ws.Range("K5").Value = 9.872375087786049
$ws.Range("M5").Value = 15.30504012954257
$ws.Range("N5").Value = 20.35834747052239
$ws.Range("O5").Value = 26.60192932555659

# Row 6
$ws.Range("B6").Value = 10.29895501941319
$ws.Range("C6").Value = 3.624500254563509
$ws.Range("D6").Value = 8.775182255342715
$ws.Range("E6").Value = 13.59202414393063
$ws.Range("F6").Value = 34.95242784061876
$ws.Range("H6").Value = 7.344005520526261
$ws.Range("J6").Value = 10.23082311650933
$ws.Range("K6").Value = 9.866439796619439
$ws.Range("M6").Value = 15.30314394833804
$ws.Range("N6").Value = 20.36107883485682
$ws.Range("O6").Value = 26.60543030325816

# Row 7
$ws.Range("B7").Value = 10.36378381545826
$ws.Range("C7").Value = 3.678277130484822
$ws.Range("D7").Value = 8.782594432335131
$ws.Range("E7").Value = 13.59151933026146
$ws.Range("F7").Value = 34.93083757674001
$ws.Range("H7").Value = 7.344005520526261
$ws.Range("J7").Value = 10.22436885113645
$ws.Range("K7").Value = 9.90772124263645
$ws.Range("M7").Value = 15.31652105182472
$ws.Range("N7").Value = 20.34227600284247
$ws.Range("O7").Value = 26.58156235483314

# Row 8
$ws.Range("B8").Value = 10.64965721480571
$ws.Range("C8").Value = 3.904976560875141
$ws.Range("D8").Value = 8.818346226665621
$ws.Range("E8").Value = 13.594279736216
$ws.Range("F8").Value = 34.85143634854319
$ws.Range("H8").Value = 7.344005520526261
$ws.Range("J8").Value = 10.19837481503335
$ws.Range("K8").Value = 10.09166828091738
$ws.Range("M8").Value = 15.38082432969454
$ws.Range("N8").Value = 20.26315009451855
$ws.Range("O8").Value = 26.48695475422834

# Row 9
$ws.Range("B9").Value = 11.2042340760616
$ws.Range("C9").Value = 4.309226284384465
$ws.Range("D9").Value = 8.899050801488965
$ws.Range("E9").Value = 13.61750108193892
$ws.Range("F9").Value = 34.75288567505749
$ws.Range("H9").Value = 7.344005520526261
$ws.Range("J9").Value = 10.15631957874163
$ws.Range("K9").Value = 10.45595706339193
$ws.Range("M9").Value = 15.52525845548937
$ws.Range("N9").Value = 20.12172366107763
$ws.Range("O9").Value = 26.34003751757035

# Row 10
$ws.Range("B10").Value = 11.60354828633762
$ws.Range("C10").Value = 4.580156980639996
$ws.Range("D10").Value = 8.964470795275796
$ws.Range("E10").Value = 13.64529332349813
$ws.Range("F10").Value = 34.71539096037259
$ws.Range("H10").Value = 7.344005520526261
$ws.Range("J10").Value = 10.13084010677876
$ws.Range("K10").Value = 10.72327486676931
$ws.Range("M10").Value = 15.64195538442655
$ws.Range("N10").Value = 20.02610902314024
$ws.Range("O10").Value = 26.25571155542312

# Row 11
$ws.Range("B11").Value = 11.78256805715219
$ws.Range("C11").Value = 4.697513649901412
$ws.Range("D11").Value = 8.995485210240734
$ws.Range("E11").Value = 13.66023501527012
$ws.Range("F11").Value = 34.70591804519928
$ws.Range("H11").Value = 7.344005520526261
$ws.Range("J11").Value = 10.12042131319568
$ws.Range("K11").Value = 10.84429882644518
$ws.Range("M11").Value = 15.69720354020909
$ws.Range("N11").Value = 19.98439380994128
$ws.Range("O11").Value = 26.22249148375044

# Row 12
$ws.Range("B12").Value = 11.84991279533192
$ws.Range("C12").Value = 4.741091074227715
$ws.Range("D12").Value = 9.00740300765702
$ws.Range("E12").Value = 13.66622043748256
$ws.Range("F12").Value = 34.70342071162597
$ws.Range("H12").Value = 7.344005520526261
$ws.Range("J12").Value = 10.1166441691064
$ws.Range("K12").Value = 10.89000250858793
$ws.Range("M12").Value = 15.71842307974651
$ws.Range("N12").Value = 19.96885202963487
$ws.Range("O12").Value = 26.21065204209316

# Row 13
$ws.Range("B13").Value = 11.83542977930702
$ws.Range("C13").Value = 4.731744501398123
$ws.Range("D13").Value = 9.004828704136619
$ws.Range("E13").Value = 13.6649168662529
$ws.Range("F13").Value = 34.70391010226213
$ws.Range("H13").Value = 7.344005520526261
$ws.Range("J13").Value = 10.11745016756517
$ws.Range("K13").Value = 10.88016563671548
$ws.Range("M13").Value = 15.71384001287359
$ws.Range("N13").Value = 19.97218791306402
$ws.Range("O13").Value = 26.21316893728537

# Row 14
$ws.Range("B14").Value = 11.78811786462007
$ws.Range("C14").Value = 4.701116167576486
$ws.Range("D14").Value = 8.996462262657367
$ws.Range("E14").Value = 13.66072089896273
$ws.Range("F14").Value = 34.70569075162078
$ws.Range("H14").Value = 7.344005520526261
$ws.Range("J14").Value = 10.12010719549661
$ws.Range("K14").Value = 10.84806164597036
$ws.Range("M14").Value = 15.69894338084731
$ws.Range("N14").Value = 19.98311007697809
$ws.Range("O14").Value = 26.22150260356738

# Row 15
$ws.Range("B15").Value = 11.75907790561429
$ws.Range("C15").Value = 4.682242617759496
$ws.Range("D15").Value = 8.991359931379378
$ws.Range("E15").Value = 13.65819327468692
$ws.Range("F15").Value = 34.70692335237781
$ws.Range("H15").Value = 7.344005520526261
$ws.Range("J15").Value = 10.12175660095194
$ws.Range("K15").Value = 10.82837946268013
$ws.Range("M15").Value = 15.68985722341056
$ws.Range("N15").Value = 19.98983337485787
$ws.Range("O15").Value = 26.22670365171996

# Row 16
$ws.Range("B16").Value = 11.5917900887576
$ws.Range("C16").Value = 4.572367480890213
$ws.Range("D16").Value = 8.962468574457791
$ws.Range("E16").Value = 13.64436283645857
$ws.Range("F16").Value = 34.71616259613894
$ws.Range("H16").Value = 7.344005520526261
$ws.Range("J16").Value = 10.13154455676866
$ws.Range("K16").Value = 10.71535036558702
$ws.Range("M16").Value = 15.63838722367275
$ws.Range("N16").Value = 20.02887093748389
$ws.Range("O16").Value = 26.25798609873939

# Row 17
$ws.Range("B17").Value = 11.48844104139108
$ws.Range("C17").Value = 4.503441585040825
$ws.Range("D17").Value = 8.94506101512656
$ws.Range("E17").Value = 13.63646488031504
$ws.Range("F17").Value = 34.72377263248607
$ws.Range("H17").Value = 7.344005520526261
$ws.Range("J17").Value = 10.13784910530725
$ws.Range("K17").Value = 10.64583171712846
$ws.Range("M17").Value = 15.60735695874919
$ws.Range("N17").Value = 20.05327439122065
$ws.Range("O17").Value = 26.27849424431249

# Row 18
$ws.Range("B18").Value = 11.42875362001471
$ws.Range("C18").Value = 4.463243489770253
$ws.Range("D18").Value = 8.935167170673949
$ws.Range("E18").Value = 13.63213877595459
$ws.Range("F18").Value = 34.72886360613176
$ws.Range("H18").Value = 7.344005520526261
$ws.Range("J18").Value = 10.14158564233416
$ws.Range("K18").Value = 10.60579424414629
$ws.Range("M18").Value = 15.58971346744304
$ws.Range("N18").Value = 20.06747825672864
$ws.Range("O18").Value = 26.29077375187989

# Row 19
$ws.Range("B19").Value = 11.4085047561176
$ws.Range("C19").Value = 4.449538530390368
$ws.Range("D19").Value = 8.931837853289322
$ws.Range("E19").Value = 13.63071132557151
$ws.Range("F19").Value = 34.73070994296002
$ws.Range("H19").Value = 7.344005520526261
$ws.Range("J19").Value = 10.14286972733365
$ws.Range("K19").Value = 10.59223058435545
$ws.Range("M19").Value = 15.58377515089674
$ws.Range("N19").Value = 20.07231626980429
$ws.Range("O19").Value = 26.2950144328109

# Row 20
$ws.Range("B20").Value = 11.49946846048273
$ws.Range("C20").Value = 4.510836298617558
$ws.Range("D20").Value = 8.946901863385076
$ws.Range("E20").Value = 13.63728323752577
$ws.Range("F20").Value = 34.72288864995502
$ws.Range("H20").Value = 7.344005520526261
$ws.Range("J20").Value = 10.13716655900162
$ws.Range("K20").Value = 10.65323780581657
$ws.Range("M20").Value = 15.61063912979167
$ws.Range("N20").Value = 20.05065925969359
$ws.Range("O20").Value = 26.27626104081024

# Row 21
$ws.Range("B21").Value = 11.80202714924935
$ws.Range("C21").Value = 4.710135985817121
$ws.Range("D21").Value = 8.998915043406647
$ws.Range("E21").Value = 13.66194449987885
$ws.Range("F21").Value = 34.70513816156965
$ws.Range("H21").Value = 7.344005520526261
$ws.Range("J21").Value = 10.11932219877159
$ws.Range("K21").Value = 10.85749509910928
$ws.Range("M21").Value = 15.70331089094401
$ws.Range("N21").Value = 19.97989506480246
$ws.Range("O21").Value = 26.21903470460352

# Row 22
$ws.Range("B22").Value = 11.99713775677318
$ws.Range("C22").Value = 4.83535416869249
$ws.Range("D22").Value = 9.033915399823355
$ws.Range("E22").Value = 13.67996828382964
$ws.Range("F22").Value = 34.6998889749016
$ws.Range("H22").Value = 7.344005520526261
$ws.Range("J22").Value = 10.10864028853749
$ws.Range("K22").Value = 10.99023930195395
$ws.Range("M22").Value = 15.7656095857916
$ws.Range("N22").Value = 19.93513159891778
$ws.Range("O22").Value = 26.18594931370729

# Row 23
$ws.Range("B23").Value = 11.89326541309083
$ws.Range("C23").Value = 4.768988169714985
$ws.Range("D23").Value = 9.015145338403226
$ws.Range("E23").Value = 13.67017537276742
$ws.Range("F23").Value = 34.70210974489494
$ws.Range("H23").Value = 7.344005520526261
$ws.Range("J23").Value = 10.11425181806751
$ws.Range("K23").Value = 10.91947332780525
$ws.Range("M23").Value = 15.73220537699979
$ws.Range("N23").Value = 19.95888719990934
$ws.Range("O23").Value = 26.2032124237133

# Row 24
$ws.Range("B24").Value = 11.49448380219065
$ws.Range("C24").Value = 4.50749492805372
$ws.Range("D24").Value = 8.946069260285556
$ws.Range("E24").Value = 13.6369125897403
$ws.Range("F24").Value = 34.72328606854558
$ws.Range("H24").Value = 7.344005520526261
$ws.Range("J24").Value = 10.13747478939428
$ws.Range("K24").Value = 10.64988972879168
$ws.Range("M24").Value = 15.60915464809829
$ws.Range("N24").Value = 20.05184101857308
$ws.Range("O24").Value = 26.27726914866719

# Row 25
$ws.Range("B25").Value = 11.05532522304556
$ws.Range("C25").Value = 4.204352550518758
$ws.Range("D25").Value = 8.876116895507719
$ws.Range("E25").Value = 13.60932442522094
$ws.Range("F25").Value = 34.77341839042822
$ws.Range("H25").Value = 7.344005520526261
$ws.Range("J25").Value = 10.16674370238556
$ws.Range("K25").Value = 10.35728600780915
$ws.Range("M25").Value = 15.48428168940537
$ws.Range("N25").Value = 20.15852114709161
$ws.Range("O25").Value = 26.37564078331069
